$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; this pushes existing rows 7-19 down to 8-20
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new weekly record
$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(7, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 44536
$ws.Cells.Item(7, 5).Value = 5
$ws.Cells.Item(7, 6).Value = 100112044
$ws.Cells.Item(7, 7).Value = "Perejil"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 125
$ws.Cells.Item(7, 11).Value = 2200
$ws.Cells.Item(7, 12).Value = 2200
$ws.Cells.Item(7, 13).Value = 2200
$ws.Cells.Item(7, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(7, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(7, 16).Value = 733
$ws.Cells.Item(7, 17).Value = 3
$ws.Cells.Item(7, 18).Value = "Hortaliza"
